# This template had its merge-field placeholders (e.g. "{Student}") split
# across several runs (and marked up with spell-check w:proofErr tags) due
# to how Word auto-split them while the template was being authored.
# The regeneration pass collapses each placeholder back into a single,
# clean run "{Name}" with no w:proofErr markers in between.

$d = $word.ActiveDocument

$placeholders = @(
    "{Student}",
    "{StudyForm}",
    "{StudyProgram}",
    "{StudyField}",
    "{Supervisor}",
    "{StudyStartDate}",
    "{StudyEndDate}",
    "{CreditsCount}",
    "{ApplicationYear}",
    "{ThesisTitle}",
    "{Opponent1WorkplaceAddress}",
    "{Opponent1PhoneNumber}",
    "{Opponent1Mail}",
    "{Opponent2}",
    "{Opponent2WorkplaceAddress}",
    "{Opponent2PhoneNumber}",
    "{Opponent2Mail}",
    "{Opponent3}",
    "{Opponent3WorkplaceAddress}",
    "{Opponent3PhoneNumber}",
    "{Opponent3Mail}"
)

foreach ($ph in $placeholders) {
    $d.Content.Find.Execute($ph, $false, $false, $false, $false, $false, $true, 1, $false, $ph, 2)
}

# "{CurrentDate}" sits right after the text "V Žiline dňa " inside the same
# paragraph and with identical run formatting, so a plain Find/Replace would
# fold it into that preceding run as well. Do the same text-normalizing
# replace first, then force the two pieces back apart into separate runs by
# toggling a character format on just the "{CurrentDate}" span (this does
# not change its visible formatting - bold is off both before and after).
$d.Content.Find.Execute("{CurrentDate}", $false, $false, $false, $false, $false, $true, 1, $false, "{CurrentDate}", 2)

$cdRange = $d.Content
$cdRange.Find.Execute("{CurrentDate}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $cdRange.Start
$end = $cdRange.End

$splitRange = $d.Range($start, $end)
$splitRange.Bold = 1
$splitRange = $d.Range($start, $end)
$splitRange.Bold = 0

Write-Output "Placeholders normalized into single runs."
